# Apply updated Schedule results/projections and refresh the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# Updated Team1Proj (C) / Team2Proj (E) values
$ws.Range("C2").Value = 136.30000000000001
$ws.Range("E2").Value = 111
$ws.Range("C3").Value = 138.4
$ws.Range("E3").Value = 122.6
$ws.Range("C4").Value = 134.19999999999999
$ws.Range("E4").Value = 145.5
$ws.Range("C5").Value = 80.2
$ws.Range("E5").Value = 106.8
$ws.Range("C6").Value = 149.4
$ws.Range("E6").Value = 144.30000000000001
$ws.Range("C7").Value = 133
$ws.Range("E7").Value = 89.5
$ws.Range("C10").Value = 133.85636363636365
$ws.Range("C11").Value = 137.30545454545452
$ws.Range("E12").Value = 132.50363636363636
$ws.Range("C13").Value = 106

# Make Schedule the active sheet/tab and update the view/selection
$ws.Activate()
$ws.Range("C8").Select()
$excel.ActiveWindow.ScrollRow = 1
